# Generate Report for Handoff
#
# The localization tool re-ran and produced a new staging file (new GUID
# base name) plus refreshed handoff/handback hash + timestamps. Update the
# Overview / zh-cn / de-de sheets to reflect the new generated report,
# keeping the existing (now-stale) hyperlink targets untouched -- only the
# cell text / hyperlink display text changes.

$oldGuid = "ee6f27c9-f274-419b-bbc6-60d5fc329341"
$newGuid = "a3da39c9-120d-42b6-8c73-c0b3bf24fdb3"

$newZhHash = "95daf71718bdddf55f9ef17add96ef8ae167d701"

$hyperlinkAddress = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/471bcbc04b5cd5a0eafb612de2957f2a34effc4e/e2e/$oldGuid.md"

# RGB(100,149,237) == #6495ED, matching the workbook's existing "HyperLink"
# cell style, expressed as the BGR long value the Font.Color property wants.
$hyperlinkColor = 15570276

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# --- Overview sheet ---
$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("G2").Value = "2016-08-16 06:51:11"
$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $hyperlinkAddress, "", "", "e2e\$newGuid.md")
$wsOverview.Range("B2").Font.Underline = 2
$wsOverview.Range("B2").Font.Color = $hyperlinkColor

# --- zh-cn sheet ---
$wsZhCn.Range("G2").Value = "$newGuid.$newZhHash.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-16 06:51:02"
$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $hyperlinkAddress, "", "", "$newGuid.md")
$wsZhCn.Range("A2").Font.Underline = 2
$wsZhCn.Range("A2").Font.Color = $hyperlinkColor

# --- de-de sheet ---
$wsDeDe.Range("G2").Value = "$newGuid.$newZhHash.de-de.xlf"
$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $hyperlinkAddress, "", "", "$newGuid.md")
$wsDeDe.Range("A2").Font.Underline = 2
$wsDeDe.Range("A2").Font.Color = $hyperlinkColor
